# Apply updated average_county_temperature (column K) values sourced from
# NOAA data, plus the two dependent cells (R32, S32) on row 32 that were
# recalculated as a result of the row-32 temperature update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value for column K (average_county_temperature)
$kUpdates = @{
    2  = 19.30324074074072
    3  = 17.25771604938272
    4  = 13.62268518518517
    11 = 13.75752314814816
    12 = 19.79629629629628
    13 = 0.8611111111111096
    15 = 0.8611111111111096
    17 = 5.486111111111112
    18 = 16.86342592592595
    19 = 16.86342592592595
    20 = 5.486111111111112
    21 = 12.41429539295394
    23 = 19.36574074074073
    29 = 12.41429539295394
    31 = 12.41429539295394
    32 = 12.41429539295394
    33 = 12.41429539295394
}

foreach ($row in $kUpdates.Keys) {
    $ws.Range("K$row").Value = $kUpdates[$row]
}

# Row 32 also has downstream recalculated values in columns R and S
$ws.Range("R32").Value = 1.466311090415359
$ws.Range("S32").Value = 1.558048815385048
